# Add two new henchmen entries (Mr. Sinister Clones / MC expansion) to the
# "data" sheet, then update the view state (frozen-pane scroll + selection)
# to match where the author was working when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Names and sets first ---------------------------------------------------------
$ws.Range("A43").Value = "Mr. Sinister Clones"
$ws.Range("E43").Value = "MC"
$ws.Range("C43").Value = 10
$ws.Range("D43").Value = 3
$ws.Range("H43").Value = 1

$ws.Range("A44").Value = "Sentinel Squad O*N*F*"
$ws.Range("E44").Value = "MC"
$ws.Range("C44").Value = 10
$ws.Range("D44").Value = 2
$ws.Range("F44").Value = 1

# --- file / loc columns filled in afterwards --------------------------------------
$ws.Range("M43").Value = "mc2,jpg"
$ws.Range("N43").Value = "6 4"

$ws.Range("M44").Value = "mc2,pjpg"
$ws.Range("N44").Value = "7 4"

# --- View state: scroll the frozen pane down and move the selection --------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("M45").Select()
